# Preliminary check-in: rename ODK settings "form_id" setting to "table_id",
# and add a "properties" sheet that will be used to generate properties.csv.

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("settings")

# --- settings sheet: rename the "form_id" setting to "table_id" -----------
# (column A holds the setting_name, column B the value, column C the
# display.title; row 2 used to be form_id / follow_map_position)
$settings.Range("A2").Value = "table_id"

# --- add a new "properties" sheet after "settings" ------------------------
$propertiesSheet = $wb.Worksheets.Add($null, $settings)
$propertiesSheet.Name = "properties"

$propertiesSheet.Range("A1").Value = "partition"
$propertiesSheet.Range("B1").Value = "aspect"
$propertiesSheet.Range("C1").Value = "key"
$propertiesSheet.Range("D1").Value = "type"
$propertiesSheet.Range("E1").Value = "value"

$propertiesSheet.Range("A2").Value = "Table"
$propertiesSheet.Range("B2").Value = "default"
$propertiesSheet.Range("C2").Value = "colOrder"
$propertiesSheet.Range("D2").Value = "array"
$propertiesSheet.Range("E2").Value = '["FMP_FOL_date","FMP_FOL_B_focal_AnimID","FMP_seq_num","FMP_xcoord","FMP_ycoord","FMP_meters_to_next_seq_num","FMP_community_id","FMP_xcoord_old","FMP_ycoord_old"]'

# make "properties" the active/selected tab, like in the edited workbook
$propertiesSheet.Activate()
